$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage
# (matches source workbook where these cells are stored as text strings,
# not numbers), without leaving a stray NumberFormat on the cell style.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "54.892.65"
$ws.Range("E2").Value = "  -2.88%  "

Set-TextValue $ws.Range("D3") "2.344.06"
$ws.Range("E3").Value = "  -5.62%  "

$ws.Range("E4").Value = "  -0.14%  "

Set-TextValue $ws.Range("D5") "470.80"
$ws.Range("E5").Value = "  -3.91%  "

Set-TextValue $ws.Range("D6") "144.18"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("E7").Value = "  +0.39%  "

Set-TextValue $ws.Range("D8") "0.501"
$ws.Range("E8").Value = "  -2.41%  "

Set-TextValue $ws.Range("D9") "2.339.67"
$ws.Range("E9").Value = "  -6.56%  "

Set-TextValue $ws.Range("D10") "0.0961"
$ws.Range("E10").Value = "  -2.11%  "

Set-TextValue $ws.Range("D11") "5.36"
$ws.Range("E11").Value = "  -7.43%  "

$ws.Range("E12").Value = "  -4.24%  "

$ws.Range("E13").Value = "  +0.46%  "

Set-TextValue $ws.Range("D14") "2.750.45"
$ws.Range("E14").Value = "  -5.61%  "

Set-TextValue $ws.Range("D15") "55.257.03"
$ws.Range("E15").Value = "  -2.20%  "

Set-TextValue $ws.Range("D16") "19.97"
$ws.Range("E16").Value = "  -6.02%  "

$ws.Range("E17").Value = "  -5.40%  "

Set-TextValue $ws.Range("D18") "2.350.61"
$ws.Range("E18").Value = "  -6.07%  "

Set-TextValue $ws.Range("D19") "4.52"
$ws.Range("E19").Value = "  -0.96%  "

Set-TextValue $ws.Range("D20") "311.35"
$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("E21").Value = "  -5.83%  "

Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  +0.13%  "

Set-TextValue $ws.Range("D23") "5.59"
$ws.Range("E23").Value = "  -3.98%  "

Set-TextValue $ws.Range("D24") "56.37"
$ws.Range("E24").Value = "  -4.27%  "

$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("E26").Value = "  -4.99%  "

Set-TextValue $ws.Range("D27") "0.152"
$ws.Range("E27").Value = "  -8.08%  "

Set-TextValue $ws.Range("D28") "2.465.95"
$ws.Range("E28").Value = "  -4.92%  "

$ws.Range("E30").Value = "  +0.15%  "

Set-TextValue $ws.Range("D31") "0.0₃0750"
$ws.Range("E31").Value = "  -5.38%  "

Set-TextValue $ws.Range("D32") "148.36"
$ws.Range("E32").Value = "  -0.62%  "

Set-TextValue $ws.Range("D33") "17.98"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("E34").Value = "  -3.40%  "

Set-TextValue $ws.Range("D35") "4.99"
$ws.Range("E35").Value = "  -4.51%  "

$ws.Range("E36").Value = "  -5.91%  "

Set-TextValue $ws.Range("D37") "3.52"
$ws.Range("E37").Value = "  -6.01%  "

Set-TextValue $ws.Range("D38") "0.813"
$ws.Range("E38").Value = "  -6.20%  "

Set-TextValue $ws.Range("D39") "0.998"
$ws.Range("E39").Value = "  +0.48%  "

Set-TextValue $ws.Range("D40") "33.33"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("E42").Value = "  -5.17%  "

$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("E44").Value = "  -5.59%  "

Set-TextValue $ws.Range("D45") "0.571"
$ws.Range("E45").Value = "  -6.91%  "

Set-TextValue $ws.Range("D46") "10.15"

Set-TextValue $ws.Range("D47") "250.38"
$ws.Range("E47").Value = "  -3.36%  "

Set-TextValue $ws.Range("D48") "0.0221"
$ws.Range("E48").Value = "  -3.60%  "

Set-TextValue $ws.Range("D49") "4.40"
$ws.Range("E49").Value = "  -9.22%  "

Set-TextValue $ws.Range("D50") "16.70"
$ws.Range("E50").Value = "  -5.21%  "

Set-TextValue $ws.Range("D51") "1.766.47"
$ws.Range("E51").Value = "  -6.86%  "
